$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Task 5: Week 4" section starting at row 28.
# Values are written in the same order the original author entered them so
# the shared-string table comes out in the same sequence.
$ws.Cells.Item(28, 1).Value = "Task 5: Week 4"
$ws.Range("A28").Font.Bold = $true

$ws.Cells.Item(29, 1).Value = "Create Icons for notification"
$ws.Cells.Item(29, 2).Value = " 1 . 45 hr"

$ws.Cells.Item(30, 1).Value = "Make sure application runs from taskbar"

$ws.Cells.Item(31, 2).Value = "3 hrs"
$ws.Cells.Item(31, 1).Value = "Test notifications"

$ws.Cells.Item(32, 1).Value = "Make necessary changes to UI"
$ws.Cells.Item(32, 2).Value = " 6 hrs"

$ws.Cells.Item(30, 2).Value = "6 hrs"

# Update the selection / scroll position to match the saved view
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("A33").Select()
